$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two more rows following the existing alternating email/password pattern
$ws.Range("A11").Value = "test@example.com"
$ws.Range("B11").Value = "password123"

$ws.Range("A12").Value = "tester@ample.com"
$ws.Range("B12").Value = "password754"
